$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 346765
$ws.Range("C2").Value = 228720
$ws.Range("D2").Value = 575485
$ws.Range("E2").Value = 39.74386821550519
$ws.Range("F2").Value = 60.25613178449481
$ws.Range("J2").Value = 1427
$ws.Range("K2").Value = 1322
$ws.Range("L2").Value = 2749

$ws.Range("B3").Value = 1000446.64
$ws.Range("C3").Value = 623805.56
$ws.Range("D3").Value = 1624252.2
$ws.Range("E3").Value = 38.40570817758474
$ws.Range("F3").Value = 61.59429182241526
$ws.Range("G3").Value = 172.7376530255334
$ws.Range("H3").Value = 188.5085403659539
$ws.Range("I3").Value = 182.2405796849614
$ws.Range("J3").Value = 17636
$ws.Range("K3").Value = 3034
$ws.Range("L3").Value = 20670

$ws.Range("B4").Value = 1743215.87
$ws.Range("C4").Value = 1090403
$ws.Range("D4").Value = 2833618.87
$ws.Range("E4").Value = 38.48093374674626
$ws.Range("F4").Value = 61.51906625325374
$ws.Range("G4").Value = 74.79853818552048
$ws.Range("H4").Value = 74.2437627657983
$ws.Range("I4").Value = 74.45682819453776
$ws.Range("J4").Value = 34974
$ws.Range("K4").Value = 3932
$ws.Range("L4").Value = 38906

$ws.Range("B5").Value = 3719900.83
$ws.Range("C5").Value = 1194732.01
$ws.Range("D5").Value = 4914632.84
$ws.Range("E5").Value = 24.30969003983622
$ws.Range("F5").Value = 75.69030996016379
$ws.Range("G5").Value = 9.567931306131761
$ws.Range("H5").Value = 113.3930108151206
$ws.Range("I5").Value = 73.44015075676002
$ws.Range("J5").Value = 50119
$ws.Range("K5").Value = 3672
$ws.Range("L5").Value = 53791

$ws.Range("B6").Value = 3877218.12
$ws.Range("C6").Value = 1324727.5
$ws.Range("D6").Value = 5201945.62
$ws.Range("E6").Value = 25.46600054615719
$ws.Range("F6").Value = 74.53399945384281
$ws.Range("G6").Value = 10.88072378675113
$ws.Range("H6").Value = 4.229072149754054
$ws.Range("I6").Value = 5.846068045237751
$ws.Range("J6").Value = 57665
$ws.Range("K6").Value = 3867
$ws.Range("L6").Value = 61532

$ws.Range("B7").Value = 3374358.05
$ws.Range("C7").Value = 1182091.18
$ws.Range("D7").Value = 4556449.23
$ws.Range("E7").Value = 25.94325362427006
$ws.Range("F7").Value = 74.05674637572996
$ws.Range("G7").Value = -10.76721967348002
$ws.Range("H7").Value = -12.96961002544784
$ws.Range("I7").Value = -12.40874928638721
$ws.Range("J7").Value = 57771
$ws.Range("K7").Value = 3793
$ws.Range("L7").Value = 61564

$ws.Range("B8").Value = 5537522.38
$ws.Range("C8").Value = 1034714.53
$ws.Range("D8").Value = 6572236.91
$ws.Range("E8").Value = 15.74371928719776
$ws.Range("F8").Value = 84.25628071280224
$ws.Range("G8").Value = -12.46745196085466
$ws.Range("H8").Value = 64.10595135273212
$ws.Range("I8").Value = 44.24031912235311
$ws.Range("J8").Value = 64032
$ws.Range("K8").Value = 3217
$ws.Range("L8").Value = 67249

$ws.Range("B9").Value = 2240499.41
$ws.Range("C9").Value = 320136.99
$ws.Range("D9").Value = 2560636.4
$ws.Range("E9").Value = 12.5022431923564
$ws.Range("F9").Value = 87.49775680764358
$ws.Range("G9").Value = -69.06035619312314
$ws.Range("H9").Value = -59.53967756244083
$ws.Range("I9").Value = -61.03858648029168
$ws.Range("J9").Value = 21853
$ws.Range("K9").Value = 922
$ws.Range("L9").Value = 22775

